# "Generate Report for Handoff" was re-run. Every file row that is still
# pending handoff (currently "Handback transform failed" / "Ready for
# handoff") gets stamped with the new, common handoff run timestamp in its
# "Latest Handoff Date(time)" column - replacing whatever distinct value
# (including any stale duplicates) it had before.

$wb = $excel.ActiveWorkbook

# Rows (1-based, header is row 1) whose handoff timestamp is refreshed by
# this handoff run - same set on every per-locale sheet.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-23 04:32:23"
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-23 04:32:18"
}

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-23 04:32:23"
}
